# Auto-generated Excel COM-interop edit script
# Applies cell value updates to match the target diff (scheduled runner price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 177.2
$ws.Range("I6").Value = 163.55556
$ws.Range("K6").Value = 490.66668
$ws.Range("M6").Value = -378.66668

$ws.Range("H17").Value = 20415.16
$ws.Range("J17").Value = 21692.723
$ws.Range("L17").Value = 65078.16900000001
$ws.Range("N17").Value = -65414.16900000001

$ws.Range("H19").Value = 1902.5
$ws.Range("I19").Value = 1409.6666
$ws.Range("K19").Value = 1409.6666
$ws.Range("M19").Value = -1234.6666

$ws.Range("H45").Value = 1894
$ws.Range("J45").Value = 1888
$ws.Range("L45").Value = 5664
$ws.Range("N45").Value = -6048

$ws.Range("H70").Value = 50611010
$ws.Range("I70").Value = 144599840
$ws.Range("J70").Value = 1638.3077
$ws.Range("K70").Value = 433799520
$ws.Range("L70").Value = 4914.9231
$ws.Range("M70").Value = -433799250
$ws.Range("N70").Value = -5454.9231

$ws.Range("H73").Value = 50611010
$ws.Range("I73").Value = 144599840
$ws.Range("J73").Value = 1638.3077
$ws.Range("K73").Value = 433799520
$ws.Range("L73").Value = 4914.9231
$ws.Range("M73").Value = -433798584
$ws.Range("N73").Value = -6786.9231

$ws.Range("H74").Value = 7597.5
$ws.Range("I74").Value = 7597.5
$ws.Range("K74").Value = 7597.5
$ws.Range("M74").Value = -6661.5

$ws.Range("H77").Value = 7597.5
$ws.Range("I77").Value = 7597.5
$ws.Range("K77").Value = 37987.5
$ws.Range("M77").Value = -33307.5

$ws.Range("H100").Value = 9876.182000000001
$ws.Range("J100").Value = 13805.429
$ws.Range("L100").Value = 13805.429
$ws.Range("N100").Value = -14887.429

$ws.Range("H132").Value = 5230.7354
$ws.Range("I132").Value = 3139.348
$ws.Range("J132").Value = 9603.637000000001
$ws.Range("K132").Value = 9418.044
$ws.Range("L132").Value = 28810.911
$ws.Range("M132").Value = -6888.044
$ws.Range("N132").Value = -33870.911

$ws.Range("H137").Value = 2421.818
$ws.Range("I137").Value = 2530.125
$ws.Range("J137").Value = 2133
$ws.Range("K137").Value = 7590.375
$ws.Range("L137").Value = 6399
$ws.Range("M137").Value = -5040.375
$ws.Range("N137").Value = -11499

$ws.Range("H141").Value = 13532659
$ws.Range("I141").Value = 18524650
$ws.Range("K141").Value = 55573950
$ws.Range("M141").Value = -55568770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3120.9
$ws.Range("I32").Value = 2743.7212
$ws.Range("J32").Value = 5677.3335
$ws.Range("K32").Value = 2743.7212
$ws.Range("L32").Value = 5677.3335
$ws.Range("M32").Value = -2456.7212
$ws.Range("N32").Value = -6251.3335

$ws.Range("H34").Value = 500
$ws.Range("I34").Value = 500
$ws.Range("K34").Value = 500
$ws.Range("M34").Value = -229

$ws.Range("H122").Value = 6573.857
$ws.Range("I122").Value = 6573.857
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19721.571
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17271.571

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0

$ws.Range("H132").Value = 4767431
$ws.Range("I132").Value = 5024.4
$ws.Range("K132").Value = 15073.2
$ws.Range("M132").Value = -12543.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 860871.9399999999
$ws.Range("I105").Value = 1837711.6
$ws.Range("J105").Value = 6137.25
$ws.Range("K105").Value = 1837711.6
$ws.Range("L105").Value = 6137.25
$ws.Range("M105").Value = -1835964.6
$ws.Range("N105").Value = -9631.25

$ws.Range("H134").Value = 4548304
$ws.Range("I134").Value = 2794.5334
$ws.Range("J134").Value = 14288681
$ws.Range("K134").Value = 8383.600199999999
$ws.Range("L134").Value = 42866043
$ws.Range("M134").Value = -5848.600199999999
$ws.Range("N134").Value = -42871113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28574592
$ws.Range("I31").Value = 50003096
$ws.Range("J31").Value = 3256.0667
$ws.Range("K31").Value = 50003096
$ws.Range("L31").Value = 3256.0667
$ws.Range("M31").Value = -50002801
$ws.Range("N31").Value = -3846.0667

$ws.Range("H34").Value = 28574592
$ws.Range("I34").Value = 50003096
$ws.Range("J34").Value = 3256.0667
$ws.Range("K34").Value = 50003096
$ws.Range("L34").Value = 3256.0667
$ws.Range("M34").Value = -50002894
$ws.Range("N34").Value = -3660.0667

$ws.Range("H58").Value = 2807.2
$ws.Range("I58").Value = 2678.8333
$ws.Range("K58").Value = 2678.8333
$ws.Range("M58").Value = -2475.8333

$ws.Range("H107").Value = 1192.72
$ws.Range("I107").Value = 877.05
$ws.Range("K107").Value = 877.05
$ws.Range("M107").Value = 1042.95

$ws.Range("H132").Value = 3225.1177
$ws.Range("I132").Value = 3166.2144
$ws.Range("K132").Value = 9498.643199999999
$ws.Range("M132").Value = -6968.643199999999

$ws.Range("H135").Value = 130000
$ws.Range("J135").Value = 130000
$ws.Range("L135").Value = 130000
$ws.Range("N135").Value = -140140

$ws.Range("H136").Value = 2807.2
$ws.Range("I136").Value = 2678.8333
$ws.Range("K136").Value = 8036.499899999999
$ws.Range("M136").Value = -5486.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1284.5
$ws.Range("I12").Value = 181.125
$ws.Range("K12").Value = 543.375
$ws.Range("M12").Value = -370.375

$ws.Range("H19").Value = 7902.3335
$ws.Range("I19").Value = 2101
$ws.Range("K19").Value = 6303
$ws.Range("M19").Value = -6129

$ws.Range("H107").Value = 8274022
$ws.Range("J107").Value = 11376646
$ws.Range("L107").Value = 34129938
$ws.Range("N107").Value = -34133778

$ws.Range("H119").Value = 997.5
$ws.Range("I119").Value = 997.5
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 2992.5
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = 1845.5

$ws.Range("H121").Value = 4169.684
$ws.Range("J121").Value = 4892.75
$ws.Range("L121").Value = 14678.25
$ws.Range("N121").Value = -17298.25

$ws.Range("H133").Value = 99999.5
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 26991.695
$ws.Range("I134").Value = 2069.2856
$ws.Range("K134").Value = 6207.8568
$ws.Range("M134").Value = -1137.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 29999
$ws.Range("J38").Value = 29999
$ws.Range("L38").Value = 29999
$ws.Range("N38").Value = -30925

$ws.Range("H102").Value = 4280.5
$ws.Range("I102").Value = 4200.6665
$ws.Range("K102").Value = 4200.6665
$ws.Range("M102").Value = -2578.6665

$ws.Range("H113").Value = 2060611.6
$ws.Range("I113").Value = 3399.4
$ws.Range("J113").Value = 4632127
$ws.Range("K113").Value = 3399.4
$ws.Range("L113").Value = 4632127
$ws.Range("M113").Value = -1229.4
$ws.Range("N113").Value = -4636467

$ws.Range("H132").Value = 12502932
$ws.Range("I132").Value = 3350.2856
$ws.Range("J132").Value = 100000000
$ws.Range("K132").Value = 10050.8568
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -7520.856800000001
$ws.Range("N132").Value = -300005060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4666.3335
$ws.Range("J22").Value = 999.5
$ws.Range("L22").Value = 999.5
$ws.Range("N22").Value = -1589.5

$ws.Range("H27").Value = 4666.3335
$ws.Range("J27").Value = 999.5
$ws.Range("L27").Value = 999.5
$ws.Range("N27").Value = -1213.5

$ws.Range("H46").Value = 3425.1428
$ws.Range("J46").Value = 3415.4
$ws.Range("L46").Value = 3415.4
$ws.Range("N46").Value = -3791.4

$ws.Range("H82").Value = 3591.9473
$ws.Range("J82").Value = 4749
$ws.Range("L82").Value = 4749
$ws.Range("N82").Value = -5471

$ws.Range("H85").Value = 3591.9473
$ws.Range("J85").Value = 4749
$ws.Range("L85").Value = 4749
$ws.Range("N85").Value = -7245

$ws.Range("H122").Value = 3566.6445
$ws.Range("I122").Value = 3273.8809
$ws.Range("K122").Value = 9821.6427
$ws.Range("M122").Value = -7371.6427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 99998.5
$ws.Range("J86").Value = 99998.5
$ws.Range("L86").Value = 99998.5
$ws.Range("N86").Value = -102244.5

$ws.Range("H89").Value = 99998.5
$ws.Range("J89").Value = 99998.5
$ws.Range("L89").Value = 499992.5
$ws.Range("N89").Value = -511224.5

$ws.Range("H113").Value = 638.8333
$ws.Range("I113").Value = 637.5
$ws.Range("J113").Value = 639.5
$ws.Range("K113").Value = 1912.5
$ws.Range("L113").Value = 1918.5
$ws.Range("M113").Value = 257.5
$ws.Range("N113").Value = -6258.5

$ws.Range("H126").Value = 3068
$ws.Range("I126").Value = 3068
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9204
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -6734

$ws.Range("H132").Value = 176386.53
$ws.Range("I132").Value = 3759.848
$ws.Range("J132").Value = 838122.2
$ws.Range("K132").Value = 11279.544
$ws.Range("L132").Value = 2514366.6
$ws.Range("M132").Value = -8749.544
$ws.Range("N132").Value = -2519426.6

$ws.Range("H136").Value = 370871.2
$ws.Range("I136").Value = 14053.577
$ws.Range("K136").Value = 42160.731
$ws.Range("M136").Value = -39610.731
